{"js": "// Remove the two \"One Sample t-test\" console-output paragraphs (Source Code\n// style) that immediately follow the `t.test(samp, mu = pop_mean, ...)` and\n// `t.test(samp, mu = 4.2, ...)` calls in the \"Type I errors\" / \"Type II\n// errors\" sections of the lab. These paragraphs start with a manual line\n// break and contain the verbatim R console output (\"    One Sample t-test\",\n// \"data:  samp\", \"t = 1.1983, df = 19, p-value = 0.2455\", ...).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// Identify every SourceCode-style output paragraph that reports the\n// one-sample t-test result for the `samp` dataset (there are similar blocks\n// for the `sleep` dataset elsewhere in the doc that must be left alone).\nconst targets = [];\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const text = paragraphs.items[i].text;\n  if (text.indexOf(\"One Sample t-test\") !== -1 && text.indexOf(\"data:  samp\") !== -1) {\n    targets.push(paragraphs.items[i]);\n  }\n}\n\n// Delete them (in reverse order is not required since we already loaded the\n// text, but it's safer in case deletion invalidates later indices).\nfor (let i = targets.length - 1; i >= 0; i--) {\n  targets[i].delete();\n}\n\nawait context.sync();\n", "ps1": "# Remove the two \"One Sample t-test\" console-output paragraphs (Source Code\n# style) that immediately follow the `t.test(samp, mu = pop_mean, ...)` and\n# `t.test(samp, mu = 4.2, ...)` calls in the \"Type I errors\" / \"Type II\n# errors\" sections of the lab. These paragraphs start with a manual line\n# break and contain the verbatim R console output (\"    One Sample t-test\",\n# \"data:  samp\", \"t = 1.1983, df = 19, p-value = 0.2455\", ...).\n\n$d = $word.ActiveDocument\n\n# Identify every SourceCode-style output paragraph that reports the\n# one-sample t-test result for the `samp` dataset (there are similar blocks\n# for the `sleep` dataset elsewhere in the doc that must be left alone).\n$targets = @()\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text\n    if (($t -like \"*One Sample t-test*\") -and ($t -like \"*data:  samp*\")) {\n        $targets += $p\n    }\n}\n\n# Delete from the last match to the first so earlier ranges stay valid as\n# later-in-document paragraphs are removed.\nfor ($i = $targets.Count - 1; $i -ge 0; $i--) {\n    $targets[$i].Range.Delete()\n}\n"}
